$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.389.55"
$ws.Range("E2").Value = "  -1.61%  "
$ws.Range("D3").Value = "1.733.93"
$ws.Range("E3").Value = "  -1.66%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.006"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.62%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "320.91"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.23%  "
$ws.Range("E6").Value = "  +0.56%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4600"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +8.15%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3514"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -3.41%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07333"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -2.05%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.53"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -2.09%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.074"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.35%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.005"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.61%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.29"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -1.92%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.897"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -2.81%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.038"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -3.44%  "
$ws.Range("D16").Value = "1.740.43"
$ws.Range("E16").Value = "  -1.79%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "90.76"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.48%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001048"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.65%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06328"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.62%  "
$ws.Range("E20").Value = "  +0.52%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.58"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -2.63%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.722"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -3.40%  "
$ws.Range("D23").Value = "27.456.87"
$ws.Range("E23").Value = "  -1.45%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.05"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.56%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.106"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "161.92"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +2.93%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.78"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.23%  "
$ws.Range("D28").Value = "1.937.26"
$ws.Range("E28").Value = "  -1.59%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "124.37"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.28%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.033"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -4.66%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.041"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -6.61%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09121"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +2.88%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.662"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.62%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.375"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -3.22%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02258"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -1.48%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "11.53"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -5.63%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05971"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -1.14%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2055"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -2.25%  "
$ws.Range("B39").Value = "TheSandbox"
$ws.Range("C39").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6213"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -1.73%  "
$ws.Range("B40").Value = "InternetComputer(DFINITY)"
$ws.Range("C40").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.858"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -2.29%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.174"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.07%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.372"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.73%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.687"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -2.39%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "12.94"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -2.62%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.689"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.21%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5774"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -1.59%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "121.56"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -1.05%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.912"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -3.74%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06824"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.10%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.105"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -6.65%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "71.04"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -3.77%  "
